# Fixes a PER (Player Efficiency Rating) data bug: the team labels in
# column B were associated with the wrong rows, and the PER values in
# column C were stale/incorrect. This re-applies the corrected team
# order and recalculated PER values for rows 2-30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = "POR"
$ws.Range("C2").Value  = 15.04666666666667

$ws.Range("B3").Value  = "NJN"
$ws.Range("C3").Value  = 3.625

$ws.Range("B4").Value  = "CLE"
$ws.Range("C4").Value  = 10.85384615384615

$ws.Range("B5").Value  = "DAL"
$ws.Range("C5").Value  = 12.68571428571429

$ws.Range("B6").Value  = "MIA"
$ws.Range("C6").Value  = 10.47058823529412

$ws.Range("B7").Value  = "SEA"
$ws.Range("C7").Value  = 13.27333333333334

$ws.Range("B8").Value  = "ATL"
$ws.Range("C8").Value  = 13.85714285714286

$ws.Range("B9").Value  = "WAS"
$ws.Range("C9").Value  = 11.96666666666667

$ws.Range("B10").Value = "MIL"
$ws.Range("C10").Value = 15.87

$ws.Range("B11").Value = "LAC"
$ws.Range("C11").Value = 13.29333333333333

$ws.Range("B12").Value = "VAN"
$ws.Range("C12").Value = 11.56428571428572

$ws.Range("B13").Value = "DET"
$ws.Range("C13").Value = 14.88235294117647

$ws.Range("B14").Value = "SAS"
$ws.Range("C14").Value = 15.21428571428571

$ws.Range("B15").Value = "ORL"
$ws.Range("C15").Value = 7.942857142857143

$ws.Range("B16").Value = "UTA"
$ws.Range("C16").Value = 12.03571428571429

$ws.Range("B17").Value = "HOU"
$ws.Range("C17").Value = 12.72142857142857

$ws.Range("B18").Value = "DEN"
$ws.Range("C18").Value = 16.7

$ws.Range("B19").Value = "LAL"
$ws.Range("C19").Value = 14.82307692307692

$ws.Range("B20").Value = "GSW"
$ws.Range("C20").Value = 12.31538461538462

$ws.Range("B21").Value = "IND"
$ws.Range("C21").Value = 12.51428571428571

$ws.Range("B22").Value = "CHI"
$ws.Range("C22").Value = 10.57333333333333

$ws.Range("B23").Value = "PHI"
$ws.Range("C23").Value = 10.96923076923077

$ws.Range("B24").Value = "CHH"
$ws.Range("C24").Value = 10.025

$ws.Range("B25").Value = "BOS"
$ws.Range("C25").Value = 11.89230769230769

$ws.Range("B26").Value = "TOR"
$ws.Range("C26").Value = 7.64

$ws.Range("B27").Value = "SAC"
$ws.Range("C27").Value = 13.24285714285714

$ws.Range("B28").Value = "PHO"
$ws.Range("C28").Value = 13.21764705882353

$ws.Range("B29").Value = "NYK"
$ws.Range("C29").Value = 14.26153846153846

$ws.Range("B30").Value = "MIN"
$ws.Range("C30").Value = 11.09230769230769
